# Update cryptos list: Price (D) and Volume(1h) (E) columns per latest scrape
# (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.480.93"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.643.84"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.68"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.55"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.70"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.80%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.339"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.109.18"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "59.414.32"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.29"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.706.69"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.11%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "339.06"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.30"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.26"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.72"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.417"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.16%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0743"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.84"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.90"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "151.28"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.68%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "286.52"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.63%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.605"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0539"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.97%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.06%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.38%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.965.72"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.25"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.03%  "
